# Applies crypto price/volume updates per commit 'Updated cryptos list on Mon Aug 21 14:36:40 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.214.47'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.680.88'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  -0.74%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5304'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.25%  '
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2685'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06307'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('E10').Value = '  -3.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07540'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '1.699.83'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.483'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5673'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008131'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').Value = '26.259.53'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.005'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.854'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.224'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1261'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.624'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06439'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.345'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.286'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.30%  '
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.659'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.010'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6114'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.416'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.180'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01618'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').Value = '1.103.51'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8671'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.71%  '
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').Value = '1.833.14'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05262'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.011'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4267'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.966'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.68%  '
